# Refreshes the scraped coinranking.com snapshot held on Sheet1.
# Column D/E cells are stored as literal text (e.g. '60.933.50' with
# two '.' thousand separators, or padded '  +0.73%  ' deltas) rather
# than numbers, so every write below goes through Set-CellText, which
# forces the Text number format for any replacement that *looks*
# numeric and then restores the cell's original Style afterwards so
# no unrelated formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($Sheet, [string]$Address, [string]$Text)

    $cell = $Sheet.Range($Address)
    # Plain (optionally signed) decimal number, e.g. '19.00', '0.998',
    # '-0.04' -- NOT the '60.933.50' / '  +0.73%  ' style text values.
    $looksNumeric = $Text -match '^\s*[-+]?\d+(\.\d+)?\s*$'

    if ($looksNumeric) {
        # Without this, Excel would silently coerce e.g. '19.00' or
        # '0.998' into the number 19 / 0.998, dropping the text
        # formatting the source data relies on.
        $originalStyle = $cell.Style
        $cell.NumberFormat = '@'
        $cell.Value = $Text
        $cell.Style = $originalStyle
    } else {
        $cell.Value = $Text
    }
}

Set-CellText $ws 'D2' '61.035.30'
Set-CellText $ws 'E2' '  +1.15%  '
Set-CellText $ws 'D3' '2.625.57'
Set-CellText $ws 'E3' '  +1.24%  '
Set-CellText $ws 'E4' '  -0.06%  '
Set-CellText $ws 'D5' '532.47'
Set-CellText $ws 'E5' '  +4.43%  '
Set-CellText $ws 'D6' '155.44'
Set-CellText $ws 'E6' '  +0.80%  '
Set-CellText $ws 'D7' '0.998'
Set-CellText $ws 'E7' '  -0.04%  '
Set-CellText $ws 'D8' '0.593'
Set-CellText $ws 'E8' '  +1.53%  '
Set-CellText $ws 'E9' '  -0.21%  '
Set-CellText $ws 'D10' '0.110'
Set-CellText $ws 'E10' '  +5.41%  '
Set-CellText $ws 'E11' '  +0.60%  '
Set-CellText $ws 'D13' '3.081.38'
Set-CellText $ws 'E13' '  +1.05%  '
Set-CellText $ws 'D14' '61.008.53'
Set-CellText $ws 'E14' '  +1.12%  '
Set-CellText $ws 'D15' '21.77'
Set-CellText $ws 'E15' '  +0.53%  '
Set-CellText $ws 'D16' '0.0000145'
Set-CellText $ws 'E16' '  +3.68%  '
Set-CellText $ws 'D17' '2.625.50'
Set-CellText $ws 'E17' '  +1.08%  '
Set-CellText $ws 'D18' '4.79'
Set-CellText $ws 'E18' '  +0.80%  '
Set-CellText $ws 'D19' '356.55'
Set-CellText $ws 'E19' '  +1.57%  '
Set-CellText $ws 'D20' '10.65'
Set-CellText $ws 'E20' '  +0.99%  '
Set-CellText $ws 'D21' '6.24'
Set-CellText $ws 'E21' '  +2.04%  '
Set-CellText $ws 'E22' '  +0.08%  '
Set-CellText $ws 'D23' '61.67'
Set-CellText $ws 'E23' '  +2.41%  '
Set-CellText $ws 'E24' '  +2.17%  '
Set-CellText $ws 'E25' '  +1.28%  '
Set-CellText $ws 'D26' '2.735.86'
Set-CellText $ws 'E26' '  +0.84%  '
Set-CellText $ws 'D27' '0.997'
Set-CellText $ws 'E27' '  -0.35%  '
Set-CellText $ws 'D28' '0.0₃0867'
Set-CellText $ws 'E28' '  +2.87%  '
Set-CellText $ws 'D29' '7.40'
Set-CellText $ws 'E29' '  +0.33%  '
Set-CellText $ws 'E30' '  -0.01%  '
Set-CellText $ws 'E31' '  +7.86%  '
Set-CellText $ws 'D32' '19.50'
Set-CellText $ws 'E32' '  +0.42%  '
Set-CellText $ws 'D33' '1.62'
Set-CellText $ws 'E33' '  +3.40%  '
Set-CellText $ws 'D34' '151.56'
Set-CellText $ws 'E34' '  -0.48%  '
Set-CellText $ws 'D35' '4.19'
Set-CellText $ws 'E35' '  +3.89%  '
Set-CellText $ws 'E36' '  +1.25%  '
Set-CellText $ws 'D37' '0.930'
Set-CellText $ws 'E37' '  +9.97%  '
Set-CellText $ws 'D38' '0.891'
Set-CellText $ws 'E38' '  +3.47%  '
Set-CellText $ws 'E39' '  +1.28%  '
Set-CellText $ws 'D40' '3.81'
Set-CellText $ws 'E40' '  +1.55%  '
Set-CellText $ws 'D41' '295.22'
Set-CellText $ws 'E41' '  -1.88%  '
Set-CellText $ws 'E42' '  +2.76%  '
Set-CellText $ws 'E43' '  +2.53%  '
Set-CellText $ws 'E44' '  +0.75%  '
Set-CellText $ws 'E45' '  +0.02%  '
Set-CellText $ws 'B46' 'RenderToken'
Set-CellText $ws 'C46' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-CellText $ws 'D46' '5.03'
Set-CellText $ws 'E46' '  +2.72%  '
Set-CellText $ws 'B47' 'EnergySwap'
Set-CellText $ws 'C47' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-CellText $ws 'D47' '19.75'
Set-CellText $ws 'E47' '  -0.19%  '
Set-CellText $ws 'B48' 'VeChain'
Set-CellText $ws 'C48' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-CellText $ws 'D48' '0.0239'
Set-CellText $ws 'E48' '  +2.65%  '
Set-CellText $ws 'D49' '10.35'
Set-CellText $ws 'E49' '  +0.48%  '
Set-CellText $ws 'B50' 'dogwifhat'
Set-CellText $ws 'C50' 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-CellText $ws 'D50' '1.86'
Set-CellText $ws 'E50' '  +4.49%  '
Set-CellText $ws 'B51' 'InjectiveProtocol'
Set-CellText $ws 'C51' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-CellText $ws 'D51' '19.00'
Set-CellText $ws 'E51' '  +5.96%  '
